$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trade rows appended (rows 5-13), matching columns A:M of the existing table.
$rows = @(
    @("BTC/USDC", 103764.81, 103866.18, 0.001, 0.1013699999999953, 0, 0.1013699999999953, 0.09769207884637897, "2025-05-16", "02:10:02", "2025-05-16", "02:12:33", 2.52413535),
    @("BTC/USDC", 103878.15, 103935.21121, 0.001, 0.05706120999999985, 0, 0.05706120999999985, 0.05493090702905264, "2025-05-16", "02:12:39", "2025-05-16", "02:15:30", 2.849939516666667),
    @("BTC/USDC", 103854.94, 103918.67704, 0.001, 0.06373703999999271, 0, 0.06373703999999271, 0.06137121642937033, "2025-05-16", "02:16:31", "2025-05-16", "02:17:42", 1.191119316666666),
    @("BTC/USDC", 103832.73, 103949.8474, 0.001, 0.1171174000000028, 0, 0.1171174000000028, 0.1127942990615799, "2025-05-16", "02:31:38", "2025-05-16", "02:32:47", 1.157217583333333),
    @("BTC/USDC", 103909.64137, 103974.21, 0.001, 0.06456863000000886, 0, 0.06456863000000886, 0.06213920974868326, "2025-05-16", "06:42:39", "2025-05-16", "07:06:21", 23.70910848333333),
    @("BTC/USDC", 103916.61, 103971.42417, 0.001, 0.05481416999999783, 0, 0.05481416999999783, 0.05274822764137305, "2025-05-16", "09:04:05", "2025-05-16", "09:05:42", 1.61504025),
    @("BTC/USDC", 103886.27174, 103500, 0.001, -0.3862717399999965, 0, -0.3862717399999965, -0.3718217369150883, "2025-05-16", "09:45:41", "2025-05-16", "10:15:26", 29.75130018333333),
    @("BTC/USDC", 103369.35434, 103431.35, 0.001, 0.0619956600000005, 0, 0.0619956600000005, 0.05997489332871894, "2025-05-16", "10:45:15", "2025-05-16", "10:54:10", 8.914118349999999),
    @("BTC/USDC", 103556.4664, 103615.03574, 0.001, 0.05856934000000183, 0, 0.05856934000000183, 0.05655787806989311, "2025-05-16", "15:01:26", "2025-05-16", "15:06:01", 4.578442166666666),
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = "'" + $row[8]
    $ws.Cells.Item($r, 10).Value = "'" + $row[9]
    $ws.Cells.Item($r, 11).Value = "'" + $row[10]
    $ws.Cells.Item($r, 12).Value = "'" + $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
}
